$wb = $excel.ActiveWorkbook

# Sheets: "Overview" (summary), "zh-cn" (detail), "de-de" (detail)
# Files that moved from "Ready for handoff" to "In Translation":
#   5f076fc6-1c2f-4b17-ad95-f7caadaea66c.md  -> row 3
#   e14f6794-ec4f-406e-b04f-cf6a852fa633.md  -> row 4
# ff65b339-b1e1-4f1d-9698-14fa462ca63d.md (row 5) keeps "Ready for handoff"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
